$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell T1, matching the style of S1 (bold header style)
$ws.Range("T1").Value = "Forma Jurídica"
$ws.Range("S1").Copy()
$ws.Range("T1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows: Forma Jurídica per company (T2:T17)
$values = @(
    "Sociedade por Quotas",
    "Sociedade por Quotas",
    "Sociedade por Quotas",
    "Sociedade por Quotas",
    "Sociedade Anónima",
    "Sociedade por Quotas",
    "Sociedade por Quotas",
    "Sociedade Anónima",
    "Sociedade por Quotas",
    "Sociedade Unipessoal",
    "Sociedade por Quotas",
    "Sociedade Unipessoal",
    "Sociedade por Quotas",
    "Sociedade por Quotas",
    "Sociedade Unipessoal",
    "Sociedade Unipessoal"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 20).Value = $values[$i]
}
